# Generate Report for handback
# A new handback round completed for the "175b4d79-..." source file: its
# Latest Handoff Datetime and Latest Handback DateTime are refreshed on
# both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-08 11:47:29"
$wsZhCn.Range("G2").Value = "2016-01-08 11:48:18"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-08 11:47:42"
$wsDeDe.Range("G2").Value = "2016-01-08 11:48:39"
